$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.907.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.521.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  +4.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.44%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +3.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.520.67'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.03%  '
$ws.Range("E10").Value = '  +4.77%  '
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.335'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.972.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.964.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("E17").Value = '  +3.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.521.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.80%  '
$ws.Range("E20").Value = '  +3.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.33%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.73%  '
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.90%  '
$ws.Range("E29").Value = '  +6.49%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.76'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.54%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("E32").Value = '  +5.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.25%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.32%  '
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("E39").Value = '  +5.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.791'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '280.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '132.27'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.601'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.28%  '
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0511'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0220'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.763.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.25%  '
